$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 370; this shifts existing rows 370-384
# down to 371-385, preserving all of their data.
$ws.Range("A370").EntireRow.Insert()

# Populate the newly inserted row 370 with the new weekly record.
$ws.Range("A370").Value = 4
$ws.Range("B370").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C370").Value = "Los Lagos"
$ws.Range("D370").Value = 45041
$ws.Range("E370").Value = 10
$ws.Range("F370").Value = 100112044
$ws.Range("G370").Value = "Perejil"
$ws.Range("H370").Value = "Sin especificar"
$ws.Range("I370").Value = "Primera"
$ws.Range("J370").Value = 160
$ws.Range("K370").Value = 5000
$ws.Range("L370").Value = 5000
$ws.Range("M370").Value = 5000
$ws.Range("N370").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O370").Value = "Región de La Araucanía"
$ws.Range("P370").Value = 2500
$ws.Range("Q370").Value = 2
$ws.Range("R370").Value = "Hortaliza"
